$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999763854153,
    0.99460662222333074,
    0.97501750906507123,
    0.96643157646890354,
    0.95479562792780048,
    0.92634700784423318,
    0.92294232539046439,
    0.91817741599738523,
    0.90642316030766135,
    0.89558854516749165,
    0.89399168031110698,
    0.89121068628612665,
    0.87991962570835325,
    0.87575132008896972,
    0.87315933551709546,
    0.87065250065518462,
    0.86694420147967888,
    0.8658352490080542,
    0.99173460486095921,
    0.96852245363004386,
    0.9621852658056228,
    0.9609207412306342,
    0.96798195819367594,
    0.95496032640693684,
    0.94850319492105273,
    0.92268950151706752,
    0.91784045076282172,
    0.89635937862712689,
    0.88108478405062307,
    0.87451363733169574,
    0.86685908426637048,
    0.86517967301558008,
    0.86465964185045796
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
